$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move A6 -> C13 and A7 -> C14 (cut the two "stretch goal" tasks out of
# column A and relocate them to the bottom of column C), preserving their
# existing cell formatting.
$ws.Range("A6").Copy($ws.Range("C13"))
$ws.Range("A7").Copy($ws.Range("C14"))
$ws.Range("A6:A7").Clear()

# Column C now holds the long "Create Temperature Gauge and Water
# interaction" label, so widen it to fit (matching column A's width).
$ws.Columns("C").ColumnWidth = $ws.Columns("A").ColumnWidth

# Leave the same active-cell/selection state recorded after the move.
$ws.Range("A10").Select() | Out-Null
